$d = $word.ActiveDocument

$para1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="01282E4F" w14:textId="29878E13" w:rsidR="00DB107E" w:rsidRDefault="001667A0"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Un </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>cañón</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> solo dispara si el otro cañón dispara primero, y se hará con el propósito de no recibir daños.</w:t></w:r></w:p>
'@

$para2PlusXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4256131E" w14:textId="793AF972" w:rsidR="001667A0" w:rsidRPr="001667A0" w:rsidRDefault="001667A0"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Las balas de </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>cañón</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> tienen un movimiento </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>parabólico</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Para que el disparo defensivo del </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>cañón</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> DD, debe estar a una distancia de no mas de 0,025d de la posición de la bala ofensiva</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>DD está dotada con sensor que detecta cualquier DO,</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>dentro del rango descrito y realiza la detonación para destruirla.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>EL rango de destrucción de DD corresponde a un circulo con centro (</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>x</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">D, </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>y</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>D)</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> y radio 0,025d</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>Por otro lado, DO está dotada con un sensor de detonación que detecta la ubicación enemiga</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>dentro a una distancia de 0,05d. El rango de destrucción de DO corresponde a todo lo que se encuentre dentro de un círculo con centro en (xO, yO) y radio 0,05d.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">El </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>cañón</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> defensivo cuenta con un infiltrado que notifica un disparo de Do y además informa los parámetros con los cuales fue realizado, pero la información llega con 2 segundes de retraso.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">EL </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>cañón</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> defensivo solo disparara si el </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>cañón</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> ofensivo dispara una bala que pueda llegar a dañar el </w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>cañón</w:t></w:r><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> defensivo</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>Me contratan para hacer que el sistema de defensa sea exitoso</w:t></w:r></w:p>
'@

$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML($para1Xml)

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML($para2PlusXml)
